$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new WR player "C.White" stat row (Week 16 logged / season sim from Week 17)
$ws.Range("A8").Value = "C.White"
$ws.Range("B8:J8").Value = 0

# Move selection to reflect the new active cell after data entry
$ws.Range("J9").Select()
